# Add a new "NBA" worksheet as the last tab in the workbook (after MMA),
# mirroring the existing SoccerPage / TeamInfo / MMA "Assertions" sheets.
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "NBA"

# Header cell, same as on the other sheets.
$ws.Range("A1").Value = "Assertions"

# NBA team names (City on first line, team name on second line).
$teams = @(
  "ATLANTA`nHAWKS",
  "BOSTON`nCELTICS",
  "BROOKLYN`nNETS",
  "CHARLOTTE`nHORNETS",
  "CHICAGO`nBULLS",
  "CLEVELAND`nCAVALIERS",
  "DALLAS`nMAVERICKS",
  "DENVER`nNUGGETS",
  "DETROIT`nPISTONS",
  "GOLDEN STATE`nWARRIORS",
  "HOUSTON`nROCKETS",
  "INDIANA`nPACERS",
  "LA`nCLIPPERS",
  "LOS ANGELES`nLAKERS",
  "MEMPHIS`nGRIZZLIES",
  "MIAMI`nHEAT",
  "MILWAUKEE`nBUCKS",
  "MINNESOTA`nTIMBERWOLVES",
  "NEW ORLEANS`nPELICANS",
  "NEW YORK`nKNICKS",
  "OKLAHOMA CITY`nTHUNDER",
  "ORLANDO`nMAGIC",
  "PHILADELPHIA`n76ERS",
  "PHOENIX`nSUNS",
  "PORTLAND`nTRAIL BLAZERS",
  "SACRAMENTO`nKINGS",
  "SAN ANTONIO`nSPURS",
  "TORONTO`nRAPTORS",
  "UTAH`nJAZZ",
  "WASHINGTON`nWIZARDS"
)

for ($i = 0; $i -lt $teams.Count; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $teams[$i]
}

$lastRow = 1 + $teams.Count
$dataRange = $ws.Range("A2:A" + $lastRow)
$dataRange.WrapText = $true
$dataRange.Font.Color = 0
$dataRange.RowHeight = 29

# Column sized to fit the wrapped team names.
$ws.Columns.Item(1).ColumnWidth = 22.09

# Match the saved view: scrolled down so row 16 is at top, with A32 selected.
[void]$ws.Range("A32").Select()
$excel.ActiveWindow.ScrollRow = 16

# Portrait page orientation, as set for the new sheet.
$ws.PageSetup.Orientation = 1
